# Update the "Förändrad" (Changed) date column (C) for all data rows (2-420)
# from 2023-09-09 (serial 45178) to 2023-09-10 (serial 45179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 420
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
